# TIYCS_battery-top-pos BOM update ("Battery PCB new gerber")
#
# Net effect (derived from the target diff):
#   - D3 (1N4148 / D_SOD-323F) row removed
#   - R1 (10K / 805) row removed
#   - C6 (C_Polarized / CP_Radial_D12.5mm_P5.00mm) row added, right after C5
#   - J4 (XT60PW-M / AMASS_XT60PW-M) row added, right after J2
#   - Several surviving rows get refreshed Mid X / Mid Y / Rotation / Val / Package
#     values coming from the regenerated pick-and-place data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Remove the two obsolete rows first (bottom-most first so row numbers of
#    rows still to be touched above stay stable).
#    Before removal, data rows are:
#      3 C1  4 C2  5 C3  6 C4  7 C5  8 D1  9 D2  10 D3
#      11 J1 12 J2 13 L1 14 R1 15 R2 16 R3 17 S1 18 U3
# ---------------------------------------------------------------------------
$ws.Rows(14).Delete()   # R1
$ws.Rows(10).Delete()   # D3

# Rows are now:
#   3 C1  4 C2  5 C3  6 C4  7 C5  8 D1  9 D2
#   10 J1 11 J2 12 L1 13 R2 14 R3 15 S1 16 U3

# ---------------------------------------------------------------------------
# 2) Insert the two new rows.
# ---------------------------------------------------------------------------
$ws.Rows(12).Insert()   # new blank row that will hold J4, right after J2 (row 11)
$ws.Rows(8).Insert()    # new blank row that will hold C6, right after C5 (row 7)

# Rows are now:
#   3 C1  4 C2  5 C3  6 C4  7 C5  8 (new C6)  9 D1  10 D2
#   11 J1 12 J2 13 (new J4) 14 L1 15 R2 16 R3 17 S1 18 U3

# ---------------------------------------------------------------------------
# 3) Write the new / refreshed row contents.
# ---------------------------------------------------------------------------

# C1 - position refreshed
$ws.Cells.Item(3, 4).Value = "20.660000"
$ws.Cells.Item(3, 5).Value = "-89.190000"

# C5 - position refreshed
$ws.Cells.Item(7, 4).Value = "69.610000"
$ws.Cells.Item(7, 5).Value = "-84.950000"

# C6 - brand new row
$ws.Cells.Item(8, 1).Value = "C6"
$ws.Cells.Item(8, 2).Value = "C_Polarized"
$ws.Cells.Item(8, 3).Value = "CP_Radial_D12.5mm_P5.00mm"
$ws.Cells.Item(8, 4).Value = "6.960000"
$ws.Cells.Item(8, 5).Value = "-36.806041"
$ws.Cells.Item(8, 6).Value = "-90.000000"
$ws.Cells.Item(8, 7).Value = "top"

# D1 - rotation flipped
$ws.Cells.Item(9, 6).Value = "90.000000"

# J1 - position + rotation refreshed
$ws.Cells.Item(11, 4).Value = "84.810000"
$ws.Cells.Item(11, 5).Value = "-91.510000"
$ws.Cells.Item(11, 6).Value = "90.000000"

# J2 - position refreshed
$ws.Cells.Item(12, 4).Value = "6.560000"
$ws.Cells.Item(12, 5).Value = "-87.220000"

# J4 - brand new row
$ws.Cells.Item(13, 1).Value = "J4"
$ws.Cells.Item(13, 2).Value = "XT60PW-M"
$ws.Cells.Item(13, 3).Value = "AMASS_XT60PW-M"
$ws.Cells.Item(13, 4).Value = "5.530000"
$ws.Cells.Item(13, 5).Value = "-66.870000"
$ws.Cells.Item(13, 6).Value = "90.000000"
$ws.Cells.Item(13, 7).Value = "top"

# R2 - value changed 10K -> 13K
$ws.Cells.Item(15, 2).Value = "13K"

# R3 - value + package changed
$ws.Cells.Item(16, 2).Value = "1.5K"
$ws.Cells.Item(16, 3).Value = "R_0805_2012Metric"

# S1 - position refreshed
$ws.Cells.Item(17, 4).Value = "6.690000"
$ws.Cells.Item(17, 5).Value = "-25.380000"
